$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 13132
$ws.Range("E2").Value = 931
$ws.Range("F2").Value = 931
$ws.Range("G2").Value = 842
$ws.Range("H2").Value = 616
$ws.Range("I2").Value = 618
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 7832
$ws.Range("L2").Value = 4523
$ws.Range("M2").Value = 3309
$ws.Range("N2").Value = 3309
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 200
$ws.Range("Q2").Value = 964
$ws.Range("R2").Value = -828
$ws.Range("S2").Value = 131
$ws.Range("T2").Value = 301
$ws.Range("U2").Value = 663
$ws.Range("V2").Value = 3228
$ws.Range("W2").Value = 7.09
$ws.Range("X2").Value = 4.69
$ws.Range("Y2").Value = 20.48
$ws.Range("Z2").Value = 8.41
$ws.Range("AA2").Value = 136.71
$ws.Range("AB2").Value = 1561.79
$ws.Range("AC2").Value = 1544
$ws.Range("AD2").Value = 25.71
$ws.Range("AE2").Value = 8273
$ws.Range("AF2").Value = 4.8
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 0.5
$ws.Range("AI2").Value = 12.95
$ws.Range("AJ2").Value = 40000000
$ws.Range("D3").Value = 15865
$ws.Range("E3").Value = 1424
$ws.Range("F3").Value = 1424
$ws.Range("G3").Value = 1373
$ws.Range("H3").Value = 1036
$ws.Range("I3").Value = 1034
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 9922
$ws.Range("L3").Value = 5666
$ws.Range("M3").Value = 4256
$ws.Range("N3").Value = 4256
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 200
$ws.Range("Q3").Value = 145
$ws.Range("R3").Value = -552
$ws.Range("S3").Value = 1237
$ws.Range("T3").Value = 389
$ws.Range("U3").Value = -244
$ws.Range("V3").Value = 4234
$ws.Range("W3").Value = 8.970000000000001
$ws.Range("X3").Value = 6.53
$ws.Range("Y3").Value = 27.34
$ws.Range("Z3").Value = 11.67
$ws.Range("AA3").Value = 133.12
$ws.Range("AB3").Value = 2034.25
$ws.Range("AC3").Value = 2586
$ws.Range("AD3").Value = 20.5
$ws.Range("AE3").Value = 10640
$ws.Range("AF3").Value = 4.98
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 0.47
$ws.Range("AI3").Value = 9.67
$ws.Range("AJ3").Value = 40000000
$ws.Range("D4").Value = 15477
$ws.Range("E4").Value = 816
$ws.Range("F4").Value = 816
$ws.Range("G4").Value = 632
$ws.Range("H4").Value = 460
$ws.Range("I4").Value = 469
$ws.Range("J4").Value = -10
$ws.Range("K4").Value = 12707
$ws.Range("L4").Value = 7071
$ws.Range("M4").Value = 5636
$ws.Range("N4").Value = 4502
$ws.Range("O4").Value = 1134
$ws.Range("P4").Value = 200
$ws.Range("Q4").Value = 569
$ws.Range("R4").Value = -503
$ws.Range("S4").Value = 421
$ws.Range("T4").Value = 431
$ws.Range("U4").Value = 138
$ws.Range("V4").Value = 5248
$ws.Range("W4").Value = 5.27
$ws.Range("X4").Value = 2.97
$ws.Range("Y4").Value = 10.72
$ws.Range("Z4").Value = 4.06
$ws.Range("AA4").Value = 125.45
$ws.Range("AB4").Value = 2208.2
$ws.Range("AC4").Value = 1174
$ws.Range("AD4").Value = 21.85
$ws.Range("AE4").Value = 11378
$ws.Range("AF4").Value = 2.25
$ws.Range("AG4").Value = 330
$ws.Range("AH4").Value = 1.29
$ws.Range("AI4").Value = 27.81
$ws.Range("AJ4").Value = 40000000
$ws.Range("D5").Value = 17113
$ws.Range("E5").Value = 565
$ws.Range("F5").Value = 565
$ws.Range("G5").Value = 719
$ws.Range("H5").Value = 461
$ws.Range("I5").Value = 433
$ws.Range("J5").Value = 28
$ws.Range("K5").Value = 11961
$ws.Range("L5").Value = 6225
$ws.Range("M5").Value = 5736
$ws.Range("N5").Value = 4757
$ws.Range("O5").Value = 979
$ws.Range("P5").Value = 200
$ws.Range("Q5").Value = 718
$ws.Range("R5").Value = -1295
$ws.Range("S5").Value = -671
$ws.Range("T5").Value = 254
$ws.Range("U5").Value = 464
$ws.Range("V5").Value = 4604
$ws.Range("W5").Value = 3.3
$ws.Range("X5").Value = 2.69
$ws.Range("Y5").Value = 9.34
$ws.Range("Z5").Value = 3.74
$ws.Range("AA5").Value = 108.53
$ws.Range("AB5").Value = 2370.44
$ws.Range("AC5").Value = 1081
$ws.Range("AD5").Value = 24.92
$ws.Range("AE5").Value = 12122
$ws.Range("AF5").Value = 2.22
$ws.Range("AG5").Value = 450
$ws.Range("AI5").Value = 40.83
$ws.Range("AJ5").Value = 40000000
$ws.Range("D6").Value = 17127
$ws.Range("E6").Value = 386
$ws.Range("F6").Value = 386
$ws.Range("G6").Value = -478
$ws.Range("H6").Value = -498
$ws.Range("I6").Value = -508
$ws.Range("K6").Value = 11707
$ws.Range("L6").Value = 6781
$ws.Range("M6").Value = 4927
$ws.Range("N6").Value = 4122
$ws.Range("P6").Value = 200
$ws.Range("Q6").Value = -214
$ws.Range("R6").Value = -251
$ws.Range("S6").Value = 198
$ws.Range("T6").Value = 317
$ws.Range("U6").Value = -531
$ws.Range("V6").Value = 5028
$ws.Range("W6").Value = 2.25
$ws.Range("X6").Value = -2.91
$ws.Range("Y6").Value = -11.43
$ws.Range("Z6").Value = -4.21
$ws.Range("AA6").Value = 137.63
$ws.Range("AB6").Value = 2095.4
$ws.Range("AC6").Value = -1269
$ws.Range("AD6").Value = -15.64
$ws.Range("AE6").Value = 10506
$ws.Range("AF6").Value = 1.89
$ws.Range("AG6").Value = 450
$ws.Range("AH6").Value = 2.27
$ws.Range("AI6").Value = -34.79
$ws.Range("AJ6").Value = 40000000
$ws.Range("D7").Value = 19201
$ws.Range("E7").Value = 829
$ws.Range("G7").Value = 566
$ws.Range("H7").Value = 425
$ws.Range("I7").Value = 438
$ws.Range("K7").Value = 12348
$ws.Range("L7").Value = 7207
$ws.Range("M7").Value = 5141
$ws.Range("N7").Value = 4364
$ws.Range("P7").Value = 200
$ws.Range("Q7").Value = 889
$ws.Range("R7").Value = -167
$ws.Range("S7").Value = -115
$ws.Range("T7").Value = 257
$ws.Range("U7").Value = 270
$ws.Range("W7").Value = 4.32
$ws.Range("X7").Value = 2.21
$ws.Range("Y7").Value = 10.33
$ws.Range("Z7").Value = 3.53
$ws.Range("AA7").Value = 140.19
$ws.Range("AC7").Value = 1096
$ws.Range("AD7").Value = 13.83
$ws.Range("AE7").Value = 11121
$ws.Range("AF7").Value = 1.36
$ws.Range("AG7").Value = 472
$ws.Range("AH7").Value = 3.12
$ws.Range("AI7").Value = 43.11
$ws.Range("D8").Value = 20491
$ws.Range("E8").Value = 1026
$ws.Range("G8").Value = 899
$ws.Range("H8").Value = 689
$ws.Range("I8").Value = 700
$ws.Range("K8").Value = 13022
$ws.Range("L8").Value = 7370
$ws.Range("M8").Value = 5652
$ws.Range("N8").Value = 4883
$ws.Range("P8").Value = 200
$ws.Range("Q8").Value = 1038
$ws.Range("R8").Value = -538
$ws.Range("S8").Value = -151
$ws.Range("T8").Value = 232
$ws.Range("U8").Value = 560
$ws.Range("W8").Value = 5
$ws.Range("X8").Value = 3.36
$ws.Range("Y8").Value = 15.15
$ws.Range("Z8").Value = 5.43
$ws.Range("AA8").Value = 130.4
$ws.Range("AC8").Value = 1751
$ws.Range("AD8").Value = 8.65
$ws.Range("AE8").Value = 12445
$ws.Range("AF8").Value = 1.22
$ws.Range("AG8").Value = 489
$ws.Range("AH8").Value = 3.23
$ws.Range("AI8").Value = 27.92
$ws.Range("D9").Value = 21612
$ws.Range("E9").Value = 1188
$ws.Range("G9").Value = 1112
$ws.Range("H9").Value = 850
$ws.Range("I9").Value = 859
$ws.Range("K9").Value = 13671
$ws.Range("L9").Value = 7326
$ws.Range("M9").Value = 6345
$ws.Range("N9").Value = 5551
$ws.Range("P9").Value = 200
$ws.Range("Q9").Value = 1236
$ws.Range("R9").Value = -534
$ws.Range("S9").Value = -280
$ws.Range("T9").Value = 313
$ws.Range("U9").Value = 870
$ws.Range("W9").Value = 5.5
$ws.Range("X9").Value = 3.94
$ws.Range("Y9").Value = 16.47
$ws.Range("Z9").Value = 6.37
$ws.Range("AA9").Value = 115.46
$ws.Range("AC9").Value = 2148
$ws.Range("AD9").Value = 7.05
$ws.Range("AE9").Value = 14147
$ws.Range("AF9").Value = 1.07
$ws.Range("AG9").Value = 488
$ws.Range("AH9").Value = 2.27
$ws.Range("AI9").Value = 22.7
